$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance (fill in green "done" cells) for the following labs:
# Row 2 (ТСИСА): labs 6,7,8,9,10 -> G2,H2,I2,J2,K2
# Row 5 (РИПС): lab 10 -> K5
# Row 6 (БЖ): labs 2,6 -> C6, G6
#
# These "done" cells use an existing cell format (fill + border, style index 3
# in the workbook) rather than a named style, so we copy that formatting from
# a cell that already carries it (e.g. B2) and then set the value to 1.

$cellsToMark = @("G2", "H2", "I2", "J2", "K2", "K5", "C6", "G6")

$formatSource = $ws.Range("B2")
$formatSource.Copy()

foreach ($addr in $cellsToMark) {
    $target = $ws.Range($addr)
    $target.PasteSpecial(-4122)  # xlPasteFormats
    $target.Value = 1
}

$excel.CutCopyMode = 0

# Update the active selection as recorded in the saved workbook.
$ws.Range("F9").Select()
